# Reorders data rows 2-26 on the active sheet according to a fixed
# permutation (new_row -> source_row taken from the original layout).
# Row 1 (header) and all columns A:T are otherwise untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 26
$firstCol = 1
$lastCol = 20

# Mapping: new row number -> old row number (1-indexed, as in the sheet)
$mapping = @{
    2  = 24
    3  = 14
    4  = 2
    5  = 3
    6  = 25
    7  = 12
    8  = 23
    9  = 17
    10 = 8
    11 = 9
    12 = 7
    13 = 5
    14 = 6
    15 = 4
    16 = 18
    17 = 11
    18 = 21
    19 = 20
    20 = 22
    21 = 26
    22 = 19
    23 = 10
    24 = 13
    25 = 15
    26 = 16
}

# Snapshot every source row's cell values before any writes, so
# overlapping reads/writes during reordering never clobber source data.
# NOTE: this runtime's `.Value` getter is unreliable here (it returns a
# stub description string instead of the cell's contents), so reads use
# `.Value2` instead; `.Value2` is fine for both get and set, and every
# column already carries a uniform number format across rows 2-26, so
# number formats never need to move between cells.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowData += , $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# Write back in the new order using the snapshot only.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $mapping[$r]
    $rowData = $snapshot[$srcRow]
    for ($i = 0; $i -lt $rowData.Length; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($r, $c).Value2 = $rowData[$i]
    }
}
